$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Additional Weapon type" (row 31) is now implemented: Missile Launcher added.
# Mark it "Yes" (matching the green fill used by other completed rows, e.g. row 27)
# and record the new header/source files.
$ws.Range("B31").Interior.Color = $ws.Range("B27").Interior.Color
$ws.Range("B31").Value = "Yes"
$ws.Range("C31").Value = "WeaponComponent.h, MissileLauncherComponent.h"
$ws.Range("D31").Value = "WeaponComponent.cpp, MissileLauncherComponent.cpp"

# Reflect the new active cell/selection left after adding the rows.
$ws.Range("D32").Select()
